# fix: rename entity file and add update fields in task spreadsheet
#
# The "numbas" task-def fields are being renamed to the new "scorm_*" naming,
# and a new trailing field ("scorm_allow_review") is appended as column W.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the three existing numbas_* headers to their scorm_* equivalents.
$ws.Range("T1").Value = "scorm_enabled"
$ws.Range("U1").Value = "scorm_time_delay_enabled"
$ws.Range("V1").Value = "scorm_attempt_limit"

# Add the new trailing header for the newly introduced field.
$ws.Range("W1").Value = "scorm_allow_review"

# Reflect the column addition in the sheet's used range / selection, matching
# a user who just added the column and left it selected.
$ws.Range("W1").Select() | Out-Null
